# "Generate Report for handback"
#
# The handback report workbook tracks, per language, the handoff/handback
# timestamps for each localized file. This run records that the
# 980d8046-22b9-43b4-98c5-ca34348e8d26 entry was handed off and handed
# back again, updating its "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) on both the zh-cn and de-de
# report sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 3 (980d8046-...) gets new handoff/handback timestamps
$zhcn.Range("D3").Value = "2016-01-08 19:59:01"
$zhcn.Range("G3").Value = "2016-01-08 19:59:42"

# de-de sheet, row 3 (980d8046-...) gets new handoff/handback timestamps
$dede.Range("D3").Value = "2016-01-08 19:59:11"
$dede.Range("G3").Value = "2016-01-08 19:59:58"
